# Update placeholder Jira IDs with the real OPQA ticket numbers for the
# "watch list scripts" test cases (TestCase_B25 .. TestCase_B44) on the
# "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$jiraIds = @(
    "OPQA-1434",
    "OPQA-1435",
    "OPQA-1436",
    "OPQA-1437",
    "OPQA-1438",
    "OPQA-1439",
    "OPQA-1440",
    "OPQA-1441",
    "OPQA-1442",
    "OPQA-1443",
    "OPQA-1444",
    "OPQA-1445",
    "OPQA-1447",
    "OPQA-1449",
    "OPQA-1450",
    "OPQA-1452",
    "OPQA-1453",
    "OPQA-1455",
    "OPQA-1456",
    "OPQA-1501"
)

$startRow = 26
for ($i = 0; $i -lt $jiraIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $jiraIds[$i]
}
